$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Competency" column (I) to the bulk-upload template, with a
# sample "Java" value in the first data row, matching the existing
# header/data formatting conventions.

# Header cell I1: text + same formatting as the other header cells (bold,
# wrap text) -- copy the format from H1 (the previous last header column).
$ws.Range("I1").Value2 = "Competency"
[void]$ws.Range("H1").Copy()
[void]$ws.Range("I1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Data cell I2: sample value, unformatted like the other data cells.
$ws.Range("I2").Value2 = "Java"

# Give the new column a sensible width (close to what Excel would pick).
$ws.Columns.Item(9).ColumnWidth = 14.89

# Reflect the column addition in the sheet's selection, like Excel does
# after inserting/filling a new column.
[void]$ws.Range("I1:I1048576").Select()
